$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 656.5417
$ws.Range("J17").Value = 656.5417
$ws.Range("L17").Value = 1969.6251
$ws.Range("N17").Value = -2305.6251
$ws.Range("H19").Value = 1680.2222
$ws.Range("I19").Value = 1593
$ws.Range("J19").Value = 1789.25
$ws.Range("K19").Value = 1593
$ws.Range("L19").Value = 1789.25
$ws.Range("M19").Value = -1418
$ws.Range("N19").Value = -2139.25
$ws.Range("H28").Value = 779.5
$ws.Range("I28").Value = 400.30435
$ws.Range("K28").Value = 400.30435
$ws.Range("M28").Value = 84.69565
$ws.Range("H40").Value = 2004.75
$ws.Range("I40").Value = 1933.3334
$ws.Range("K40").Value = 1933.3334
$ws.Range("M40").Value = -1758.3334
$ws.Range("H54").Value = 19999.5
$ws.Range("I54").Value = 19999.5
$ws.Range("K54").Value = 19999.5
$ws.Range("M54").Value = -19513.5
$ws.Range("H55").Value = 360
$ws.Range("I55").Value = 146.66667
$ws.Range("K55").Value = 146.66667
$ws.Range("M55").Value = 67.33332999999999
$ws.Range("H98").Value = 2353.2942
$ws.Range("I98").Value = 2318.375
$ws.Range("K98").Value = 2318.375
$ws.Range("M98").Value = -820.375
$ws.Range("H113").Value = 7495
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 7495
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 7495
$ws.Range("M113").Value = ""
$ws.Range("N113").Value = -14003
$ws.Range("H121").Value = 999
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").Value = ""
$ws.Range("H122").Value = 2353.2942
$ws.Range("I122").Value = 2318.375
$ws.Range("K122").Value = 6955.125
$ws.Range("M122").Value = -4505.125
$ws.Range("H131").Value = 3610.2856
$ws.Range("J131").Value = 6415.8335
$ws.Range("L131").Value = 19247.5005
$ws.Range("N131").Value = -29327.5005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").Value = ""
$ws.Range("H32").Value = 22978.277
$ws.Range("I32").Value = 23434.348
$ws.Range("K32").Value = 23434.348
$ws.Range("M32").Value = -23147.348
$ws.Range("H41").Value = 1856.8334
$ws.Range("I41").Value = 1478.2
$ws.Range("J41").Value = 3750
$ws.Range("K41").Value = 1478.2
$ws.Range("L41").Value = 3750
$ws.Range("M41").Value = -1064.2
$ws.Range("N41").Value = -4578
$ws.Range("H45").Value = 3233.7058
$ws.Range("I45").Value = 1576.1111
$ws.Range("K45").Value = 1576.1111
$ws.Range("M45").Value = -1199.1111
$ws.Range("H61").Value = 3947.6216
$ws.Range("I61").Value = 833.2593000000001
$ws.Range("J61").Value = 12356.4
$ws.Range("K61").Value = 833.2593000000001
$ws.Range("L61").Value = 12356.4
$ws.Range("M61").Value = -621.2593000000001
$ws.Range("N61").Value = -12780.4
$ws.Range("H88").Value = 8158.875
$ws.Range("I88").Value = 2290.6667
$ws.Range("J88").Value = 11679.8
$ws.Range("K88").Value = 2290.6667
$ws.Range("L88").Value = 11679.8
$ws.Range("M88").Value = -1884.6667
$ws.Range("N88").Value = -12491.8
$ws.Range("H91").Value = 8158.875
$ws.Range("I91").Value = 2290.6667
$ws.Range("J91").Value = 11679.8
$ws.Range("K91").Value = 2290.6667
$ws.Range("L91").Value = 11679.8
$ws.Range("M91").Value = -886.6667000000002
$ws.Range("N91").Value = -14487.8
$ws.Range("H110").Value = 626.9167
$ws.Range("I110").Value = 620.2727
$ws.Range("K110").Value = 620.2727
$ws.Range("M110").Value = 1424.7273
$ws.Range("H122").Value = 3178.647
$ws.Range("I122").Value = 2989.08
$ws.Range("K122").Value = 8967.24
$ws.Range("M122").Value = -6517.24
$ws.Range("H136").Value = 3947.6216
$ws.Range("I136").Value = 833.2593000000001
$ws.Range("J136").Value = 12356.4
$ws.Range("K136").Value = 2499.7779
$ws.Range("L136").Value = 37069.2
$ws.Range("M136").Value = 50.22209999999995
$ws.Range("N136").Value = -42169.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1305.6
$ws.Range("I80").Value = 1336.1666
$ws.Range("J80").Value = 1259.75
$ws.Range("K80").Value = 1336.1666
$ws.Range("L80").Value = 1259.75
$ws.Range("M80").Value = -338.1666
$ws.Range("N80").Value = -3255.75
$ws.Range("H83").Value = 1305.6
$ws.Range("I83").Value = 1336.1666
$ws.Range("J83").Value = 1259.75
$ws.Range("K83").Value = 6680.833000000001
$ws.Range("L83").Value = 6298.75
$ws.Range("M83").Value = -1688.833000000001
$ws.Range("N83").Value = -16282.75
$ws.Range("H86").Value = 1932
$ws.Range("I86").Value = 1712.3334
$ws.Range("K86").Value = 1712.3334
$ws.Range("M86").Value = -589.3334
$ws.Range("H89").Value = 1932
$ws.Range("I89").Value = 1712.3334
$ws.Range("K89").Value = 8561.666999999999
$ws.Range("M89").Value = -2945.666999999999
$ws.Range("H105").Value = 2585.7917
$ws.Range("I105").Value = 1564.0625
$ws.Range("K105").Value = 1564.0625
$ws.Range("M105").Value = 182.9375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 12925.069
$ws.Range("I58").Value = 1296.6818
$ws.Range("K58").Value = 1296.6818
$ws.Range("M58").Value = -1093.6818
$ws.Range("H122").Value = 2167.0667
$ws.Range("I122").Value = 2126.5454
$ws.Range("J122").Value = 2278.5
$ws.Range("K122").Value = 6379.6362
$ws.Range("L122").Value = 6835.5
$ws.Range("M122").Value = -3929.6362
$ws.Range("N122").Value = -11735.5
$ws.Range("H136").Value = 12925.069
$ws.Range("I136").Value = 1296.6818
$ws.Range("K136").Value = 3890.0454
$ws.Range("M136").Value = -1340.0454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 43333
$ws.Range("J32").Value = 43333
$ws.Range("L32").Value = 129999
$ws.Range("N32").Value = -130565
$ws.Range("H98").Value = 1537.091
$ws.Range("J98").Value = 1537.091
$ws.Range("L98").Value = 4611.272999999999
$ws.Range("N98").Value = -7607.272999999999
$ws.Range("H103").Value = 945.3333
$ws.Range("I103").Value = 656.25
$ws.Range("J103").Value = 1176.6
$ws.Range("K103").Value = 1968.75
$ws.Range("L103").Value = 3529.8
$ws.Range("M103").Value = -1089.75
$ws.Range("N103").Value = -5287.799999999999
$ws.Range("H131").Value = 3335
$ws.Range("I131").Value = 3459.6155
$ws.Range("J131").Value = 2525
$ws.Range("K131").Value = 10378.8465
$ws.Range("L131").Value = 7575
$ws.Range("M131").Value = -5338.8465
$ws.Range("N131").Value = -17655
$ws.Range("H132").Value = 1025.3636
$ws.Range("J132").Value = 1266
$ws.Range("L132").Value = 11394
$ws.Range("N132").Value = -16454

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 5362.5
$ws.Range("I55").Value = 3340
$ws.Range("J55").Value = 8733.333000000001
$ws.Range("K55").Value = 3340
$ws.Range("L55").Value = 8733.333000000001
$ws.Range("M55").Value = -3013
$ws.Range("N55").Value = -9387.333000000001
$ws.Range("H70").Value = 6225.5
$ws.Range("J70").Value = 6463.4287
$ws.Range("L70").Value = 6463.4287
$ws.Range("N70").Value = -7003.4287
$ws.Range("H73").Value = 6225.5
$ws.Range("J73").Value = 6463.4287
$ws.Range("L73").Value = 6463.4287
$ws.Range("N73").Value = -8335.4287
$ws.Range("H102").Value = 2356.5925
$ws.Range("I102").Value = 1465.9375
$ws.Range("K102").Value = 1465.9375
$ws.Range("M102").Value = 156.0625
$ws.Range("H122").Value = 2834.5312
$ws.Range("I122").Value = 2502.96
$ws.Range("K122").Value = 7508.88
$ws.Range("M122").Value = -5058.88
$ws.Range("H126").Value = 2781.3333
$ws.Range("I126").Value = 1861.7142
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 5585.142599999999
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -3115.142599999999
$ws.Range("N126").Value = -22940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 15003.5
$ws.Range("J23").Value = 15003.5
$ws.Range("L23").Value = 15003.5
$ws.Range("N23").Value = -15463.5
$ws.Range("H31").Value = 1188.3334
$ws.Range("I31").Value = 1015
$ws.Range("J31").Value = 1275
$ws.Range("K31").Value = 1015
$ws.Range("L31").Value = 1275
$ws.Range("M31").Value = -767
$ws.Range("N31").Value = -1771
$ws.Range("H105").Value = 49999
$ws.Range("J105").Value = 49999
$ws.Range("L105").Value = 49999
$ws.Range("N105").Value = -56987
$ws.Range("H122").Value = 4001
$ws.Range("I122").Value = 4000
$ws.Range("J122").Value = 4002.5
$ws.Range("K122").Value = 12000
$ws.Range("L122").Value = 12007.5
$ws.Range("M122").Value = -9550
$ws.Range("N122").Value = -16907.5
$ws.Range("H136").Value = 4212.2144
$ws.Range("I136").Value = 3923.7058
$ws.Range("K136").Value = 11771.1174
$ws.Range("M136").Value = -9221.117400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 47510
$ws.Range("I51").Value = 26887.5
$ws.Range("K51").Value = 26887.5
$ws.Range("M51").Value = -26377.5
$ws.Range("H62").Value = 4289.636
$ws.Range("I62").Value = 3547.75
$ws.Range("K62").Value = 3547.75
$ws.Range("M62").Value = -2923.75
$ws.Range("H65").Value = 4289.636
$ws.Range("I65").Value = 3547.75
$ws.Range("K65").Value = 17738.75
$ws.Range("M65").Value = -14618.75
$ws.Range("H122").Value = 58764.24
$ws.Range("I122").Value = 85096.12
$ws.Range("J122").Value = 2809
$ws.Range("K122").Value = 255288.36
$ws.Range("L122").Value = 8427
$ws.Range("M122").Value = -252838.36
$ws.Range("N122").Value = -13327
$ws.Range("H132").Value = 843.9474
$ws.Range("I132").Value = 693.45715
$ws.Range("K132").Value = 2080.37145
$ws.Range("M132").Value = 449.6285500000004
